# Natmi following Dr Hou advice
# Re-computed NATMI Cytl1->Ccr2 edge-weight table after adding the "FAPs"
# cluster to the analysis. Every row's ligand/receptor specificity numbers
# shift (because they're computed relative to ALL clusters), the old
# target-cluster ordering (ECs, M1, M2, sCs) gets a new "FAPs" row inserted
# alphabetically between "ECs" and "M1", and the former M1/M2/sCs rows slide
# down one row with freshly recomputed values. A brand new row 6 (sCs) is
# appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ RowNum=2; D="ECs";  E=3; F=1; G=1.333368333333333; H=4.000105; I=1; J=1; K=3; L=1;
       M=0.06930366666666667; N=0.207911; O=0.0001844138843618826; P=0.0001844379559408913;
       Q=0.09240731451722223; R=0.831665830655; S=0.0001844138843618826; T=0.0001844379559408913 },
    @{ RowNum=3; D="FAPs"; E=3; F=1; G=1.333368333333333; H=4.000105; I=1; J=1; K=2; L=0.6666666666666666;
       M=0.311428; N=0.934284; O=0.000828695651202472; P=0.000828803821001677;
       Q=0.4152482333133334; R=3.737234099820001; S=0.000828695651202472; T=0.000828803821001677 },
    @{ RowNum=4; D="M1";   E=3; F=1; G=1.333368333333333; H=4.000105; I=1; J=1; K=3; L=1;
       M=168.931335; N=506.794005; O=0.4495185489626108; P=0.4495772247033483;
       Q=225.247692596725; R=2027.229233370525; S=0.4495185489626108; T=0.4495772247033483 },
    @{ RowNum=5; D="M2";   E=3; F=1; G=1.333368333333333; H=4.000105; I=1; J=1; K=3; L=1;
       M=206.345828; N=619.0374839999999; O=0.5490768020453307; P=0.5491484731435671;
       Q=275.1349927706467; R=2476.21493493582; S=0.5490768020453307; T=0.5491484731435671 },
    @{ RowNum=6; D="sCs";  E=3; F=1; G=1.333368333333333; H=4.000105; I=1; J=1; K=2; L=1;
       M=0.1471425; N=0.294285; O=0.0003915394564941487; P=0.000261060376142028;
       Q=0.1961951499875; R=1.177170899925; S=0.0003915394564941487; T=0.000261060376142028 }
)

$numCols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $rows) {
    $r = $row.RowNum
    $ws.Range("A$r").Value = "ECs"
    $ws.Range("B$r").Value = "Cytl1"
    $ws.Range("C$r").Value = "Ccr2"
    $ws.Range("D$r").Value = $row.D
    foreach ($c in $numCols) {
        $ws.Range("$c$r").Value = $row[$c]
    }
}
